$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 667.42
$ws.Range("I15").Value = 667.42
$ws.Range("K15").Value = 2002.26
$ws.Range("M15").Value = -1833.26
$ws.Range("H19").Value = 3198
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 3198
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 3198
$ws.Range("N19").Value = -3548
$ws.Range("M19").ClearContents()
$ws.Range("H21").Value = 4304.2
$ws.Range("I21").Value = 4004.6667
$ws.Range("K21").Value = 4004.6667
$ws.Range("M21").Value = -3536.6667
$ws.Range("H23").Value = 4304.2
$ws.Range("I23").Value = 4004.6667
$ws.Range("K23").Value = 4004.6667
$ws.Range("M23").Value = -3770.6667
$ws.Range("H64").Value = 7511.364
$ws.Range("I64").Value = 6218.3335
$ws.Range("J64").Value = 7996.25
$ws.Range("K64").Value = 6218.3335
$ws.Range("L64").Value = 7996.25
$ws.Range("M64").Value = -5970.3335
$ws.Range("N64").Value = -8492.25
$ws.Range("H67").Value = 7511.364
$ws.Range("I67").Value = 6218.3335
$ws.Range("J67").Value = 7996.25
$ws.Range("K67").Value = 6218.3335
$ws.Range("L67").Value = 7996.25
$ws.Range("M67").Value = -5360.3335
$ws.Range("N67").Value = -9712.25
$ws.Range("H100").Value = 2528.2222
$ws.Range("I100").Value = 1917.5
$ws.Range("J100").Value = 3749.6667
$ws.Range("K100").Value = 1917.5
$ws.Range("L100").Value = 3749.6667
$ws.Range("M100").Value = -1376.5
$ws.Range("N100").Value = -4831.6667
$ws.Range("H116").Value = 7149.6665
$ws.Range("I116").Value = 3996
$ws.Range("J116").Value = 9402.286
$ws.Range("K116").Value = 3996
$ws.Range("L116").Value = 9402.286
$ws.Range("M116").Value = -554
$ws.Range("N116").Value = -16286.286
$ws.Range("H125").Value = 9807205
$ws.Range("I125").Value = 1648.375
$ws.Range("K125").Value = 14835.375
$ws.Range("M125").Value = -12375.375
$ws.Range("H132").Value = 125007464
$ws.Range("I132").Value = 125007464
$ws.Range("K132").Value = 375022392
$ws.Range("M132").Value = -375019862
$ws.Range("H135").Value = 1760.1111
$ws.Range("I135").Value = 845.1539
$ws.Range("J135").Value = 4139
$ws.Range("K135").Value = 7606.3851
$ws.Range("L135").Value = 37251
$ws.Range("M135").Value = -5071.3851
$ws.Range("N135").Value = -42321

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 28998.5
$ws.Range("J103").Value = 28998.5
$ws.Range("L103").Value = 28998.5
$ws.Range("N103").Value = -31342.5
$ws.Range("H107").Value = 2553313.2
$ws.Range("I107").Value = 3403100
$ws.Range("J107").Value = 3953
$ws.Range("K107").Value = 3403100
$ws.Range("L107").Value = 3953
$ws.Range("M107").Value = -3401180
$ws.Range("N107").Value = -7793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1934.4615
$ws.Range("I16").Value = 1666.1666
$ws.Range("K16").Value = 1666.1666
$ws.Range("M16").Value = -1379.1666
$ws.Range("H113").Value = 1934.4615
$ws.Range("I113").Value = 1666.1666
$ws.Range("K113").Value = 1666.1666
$ws.Range("M113").Value = 503.8334
$ws.Range("H132").Value = 56163.75
$ws.Range("I132").Value = 46722.61
$ws.Range("K132").Value = 140167.83
$ws.Range("M132").Value = -137637.83
$ws.Range("H134").Value = 3487.2415
$ws.Range("I134").Value = 2296.5264
$ws.Range("J134").Value = 5749.6
$ws.Range("K134").Value = 6889.5792
$ws.Range("L134").Value = 17248.8
$ws.Range("M134").Value = -4354.5792
$ws.Range("N134").Value = -22318.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 250348
$ws.Range("I9").Value = 500300
$ws.Range("J9").Value = 396
$ws.Range("K9").Value = 1500900
$ws.Range("L9").Value = 1188
$ws.Range("M9").Value = -1500676
$ws.Range("N9").Value = -1636
$ws.Range("H122").Value = 1146.2142
$ws.Range("I122").Value = 1237.3334
$ws.Range("K122").Value = 11136.0006
$ws.Range("M122").Value = -8686.000599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10531011
$ws.Range("I70").Value = 12503981
$ws.Range("K70").Value = 12503981
$ws.Range("M70").Value = -12503711
$ws.Range("H73").Value = 10531011
$ws.Range("I73").Value = 12503981
$ws.Range("K73").Value = 12503981
$ws.Range("M73").Value = -12503045
$ws.Range("H80").Value = 32782480
$ws.Range("I80").Value = 52450196
$ws.Range("K80").Value = 52450196
$ws.Range("M80").Value = -52449198
$ws.Range("H83").Value = 32782480
$ws.Range("I83").Value = 52450196
$ws.Range("K83").Value = 262250980
$ws.Range("M83").Value = -262245988
$ws.Range("H95").Value = 41414.332
$ws.Range("J95").Value = 41414.332
$ws.Range("L95").Value = 41414.332
$ws.Range("N95").Value = -46906.332
$ws.Range("H97").Value = 1083124.9
$ws.Range("I97").Value = 1832271.9
$ws.Range("J97").Value = 1023.7778
$ws.Range("K97").Value = 1832271.9
$ws.Range("L97").Value = 1023.7778
$ws.Range("M97").Value = -1831775.9
$ws.Range("N97").Value = -2015.7778
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7568.4165
$ws.Range("I40").Value = 5646.4443
$ws.Range("J40").Value = 13334.333
$ws.Range("K40").Value = 5646.4443
$ws.Range("L40").Value = 13334.333
$ws.Range("M40").Value = -5510.4443
$ws.Range("N40").Value = -13606.333
$ws.Range("H61").Value = 8549240
$ws.Range("I61").Value = 10102138
$ws.Range("J61").Value = 8301.5
$ws.Range("K61").Value = 10102138
$ws.Range("L61").Value = 8301.5
$ws.Range("M61").Value = -10101936
$ws.Range("N61").Value = -8705.5
$ws.Range("H100").Value = 5190.1816
$ws.Range("I100").Value = 5433.731
$ws.Range("J100").Value = 4285.5713
$ws.Range("K100").Value = 5433.731
$ws.Range("L100").Value = 4285.5713
$ws.Range("M100").Value = -4892.731
$ws.Range("N100").Value = -5367.5713
$ws.Range("H113").Value = 8549240
$ws.Range("I113").Value = 10102138
$ws.Range("J113").Value = 8301.5
$ws.Range("K113").Value = 10102138
$ws.Range("L113").Value = 8301.5
$ws.Range("M113").Value = -10099968
$ws.Range("N113").Value = -12641.5
$ws.Range("H135").Value = 41285.4
$ws.Range("J135").Value = 41285.4
$ws.Range("L135").Value = 41285.4
$ws.Range("N135").Value = -51425.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8437.25
$ws.Range("I62").Value = 4002
$ws.Range("J62").Value = 8550.975
$ws.Range("K62").Value = 4002
$ws.Range("L62").Value = 8550.975
$ws.Range("M62").Value = -3378
$ws.Range("N62").Value = -9798.975
$ws.Range("H65").Value = 8437.25
$ws.Range("I65").Value = 4002
$ws.Range("J65").Value = 8550.975
$ws.Range("K65").Value = 20010
$ws.Range("L65").Value = 42754.875
$ws.Range("M65").Value = -16890
$ws.Range("N65").Value = -48994.875
$ws.Range("H99").Value = 25000
$ws.Range("I99").Value = 25000
$ws.Range("K99").Value = 25000
$ws.Range("M99").Value = -22005
$ws.Range("H113").Value = 1115.5
$ws.Range("J113").Value = 1260.85
$ws.Range("L113").Value = 3782.55
$ws.Range("N113").Value = -8122.549999999999
$ws.Range("H136").Value = 5790.579
$ws.Range("I136").Value = 6144.0835
$ws.Range("J136").Value = 5184.5713
$ws.Range("K136").Value = 18432.2505
$ws.Range("L136").Value = 15553.7139
$ws.Range("M136").Value = -15882.2505
$ws.Range("N136").Value = -20653.7139
$ws.Range("H137").Value = 76666.336
$ws.Range("J137").Value = 76666.336
$ws.Range("L137").Value = 76666.336
$ws.Range("N137").Value = -86866.336
